$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, bypassing Excel's automatic
# number coercion (needed for price strings like "212.92" that must
# stay text, matching the original inline-string cell content) while
# avoiding any NumberFormat/style change on the cell.
$xlPasteValues = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues
function Set-TextValue($addr, $val) {
    $escaped = $val.Replace('"', '""')
    $ws.Range("ZZ1").Formula = '="' + $escaped + '"'
    $ws.Range("ZZ1").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial($xlPasteValues) | Out-Null
}

# --- Plain text/percentage updates (safe via direct .Value assignment) ---
$ws.Range("D2").Value = '28.603.22'
$ws.Range("E2").Value = '  +1.33%  '
$ws.Range("D3").Value = '1.572.71'
$ws.Range("E3").Value = '  -1.15%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("E5").Value = '  -0.40%  '
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +3.89%  '
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("E10").Value = '  -1.64%  '
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("D13").Value = '1.798.11'
$ws.Range("E13").Value = '  -1.09%  '
$ws.Range("D14").Value = '1.569.39'
$ws.Range("E14").Value = '  -1.39%  '
$ws.Range("E15").Value = '  -1.91%  '
$ws.Range("D16").Value = '28.584.99'
$ws.Range("E16").Value = '  +1.05%  '
$ws.Range("E17").Value = '  -1.89%  '
$ws.Range("E18").Value = '  -1.34%  '
$ws.Range("E19").Value = '  +1.38%  '
$ws.Range("E20").Value = '  -1.55%  '
$ws.Range("D21").Value = '0.0₃0692'
$ws.Range("E21").Value = '  -2.52%  '
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("E23").Value = '  -5.15%  '
$ws.Range("E24").Value = '  -2.36%  '
$ws.Range("E25").Value = '  +10.56%  '
$ws.Range("E26").Value = '  -0.31%  '
$ws.Range("E27").Value = '  -1.38%  '
$ws.Range("E28").Value = '  -2.25%  '
$ws.Range("E29").Value = '  -3.19%  '
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("E32").Value = '  -2.61%  '
$ws.Range("E33").Value = '  -0.88%  '
$ws.Range("E34").Value = '  -1.97%  '
$ws.Range("D35").Value = '1.392.47'
$ws.Range("E35").Value = '  -0.78%  '
$ws.Range("E36").Value = '  +1.39%  '
$ws.Range("E37").Value = '  -3.30%  '
$ws.Range("E38").Value = '  +0.84%  '
$ws.Range("E40").Value = '  -0.37%  '
$ws.Range("E41").Value = '  -3.09%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("E43").Value = '  -2.94%  '
$ws.Range("E44").Value = '  +0.60%  '
$ws.Range("E45").Value = '  +1.88%  '
$ws.Range("E46").Value = '  -2.67%  '
$ws.Range("E47").Value = '  -1.73%  '
$ws.Range("E48").Value = '  -1.97%  '
$ws.Range("D49").Value = '1.710.54'
$ws.Range("E49").Value = '  -1.17%  '
$ws.Range("E50").Value = '  -1.16%  '
$ws.Range("E51").Value = '  -1.00%  '

# --- Numeric-looking price updates (must remain text, use helper) ---
Set-TextValue 'D5' '212.92'
Set-TextValue 'D8' '45.62'
Set-TextValue 'D9' '24.21'
Set-TextValue 'D18' '62.33'
Set-TextValue 'D19' '230.85'
Set-TextValue 'D23' '3.90'
Set-TextValue 'D24' '9.12'
Set-TextValue 'D26' '151.46'
Set-TextValue 'D27' '15.00'
Set-TextValue 'D38' '2.36'
Set-TextValue 'D41' '0.524'
Set-TextValue 'D43' '0.791'
Set-TextValue 'D45' '0.0468'
Set-TextValue 'D48' '63.07'
Set-TextValue 'D50' '86.61'

# Clean up helper cell
$ws.Range("ZZ1").ClearContents() | Out-Null
$excel.CutCopyMode = 0
